$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d6 = @'
id="Par1">Appropriate critical care delivery for Coronavirus disease 2019 (COVID-19) is a cornerstone in saving lives.

 Earlier publications worldwide demonstrate higher mortality among patients receiving mechanical ventilation in intensive care units during “surges” in the number of cases.

 In contrast, lower mortality outcomes are evident in Japan using CRISIS [CRoss Icu Searchable Information System] data by the national registry, Japan ECMOnet for COVID-19. This highlights the need for scientific analysis of the medical factors contributing to high survival rates and social factors associated with low case “surges,” to gain insight into protective strategies for possible coming waves in the COVID-19 pandemic.


'@
$e6 = @'
[Ryuzo%Abe%NULL%1,  Naofumi%Bunya%NULL%1,  Tomoyuki%Endo%NULL%1,  Yuji%Fujino%NULL%1,  Kensuke%Fujita%NULL%1,  Kenji%Fujizuka%NULL%1,  Yoshihiro%Hagiwara%NULL%1,  Jun%Hamaguchi%NULL%1,  Yoshitaka%Hara%NULL%1,  Eiji%Hashiba%NULL%1,  Satoru%Hashimoto%NULL%1,  Noriyuki%Hattori%NULL%1,  Kota%Hoshino%NULL%1,  Shinichi%Ijuin%NULL%1,  Takanari%Ikeyama%NULL%1,  Shingo%Ichiba%NULL%1,  Wataru%Iwanaga%NULL%1,  Yoshiaki%Iwashita%NULL%1,  Masafumi%Kanamoto%NULL%1,  Hitoshi%Kaneko%NULL%1,  Kaneyuki%Kawamae%NULL%1,  Toru%Kotani%NULL%1,  Yasuaki%Koyama%NULL%1,  Keibun%Liu%NULL%1,  Tomohiko%Masuno%NULL%1,  Naoto%Morimura%NULL%1,  Tomoyuki%Nakamura%NULL%1,  Masaki%Nakane%NULL%1,  Michitaka%Nasu%NULL%1,  Osamu%Nishida%NULL%1,  Masaji%Nishimura%NULL%1,  Kanae%Ochiai%NULL%1,  Takayuki%Ogura%NULL%1,  Shinichiro%Ohshimo%NULL%1,  Keisuke%Oyama%NULL%1,  Junichi%Sasaki%NULL%1,  Ryutaro%Seo%NULL%1,  Takeshi%Shimazu%NULL%1,  Nobuaki%Shime%NULL%4,  Keiki%Shimizu%NULL%1,  Hiroyuki%Suzuki%NULL%1,  Shuhei%Takauji%NULL%1,  Shinhiro%Takeda%NULL%1,  Ichiro%Takeuchi%NULL%2,  Mumon%Takita%NULL%1,  Hayato%Taniguchi%NULL%1,  Nobuaki%Shime%shime@koto.kpu-m.ac.jp%0,  Nobuaki%Shime%shime@koto.kpu-m.ac.jp%0,  Nobuaki%Shime%shime@koto.kpu-m.ac.jp%0]
'@
$ws.Range("D6").Value = $d6
$ws.Range("E6").Value = $e6

$d12 = @'
Background
id="Par1">When severe, COVID-19 shares many clinical features with bacterial sepsis.

 Yet, secondary bacterial infection is uncommon.

 However, as epithelium is injured and barrier function is lost, bacterial products entering the circulation might contribute to the pathophysiology of COVID-19.
Methods
id="Par2">We studied 19 adults, severely ill patients with COVID-19 infection, who were admitted to King Chulalongkorn Memorial Hospital, Bangkok, Thailand, between 13th March and 17th April 2020. Blood samples on days 1, 3, and 7 of enrollment were analyzed for endotoxin activity assay (EAA), (1 → 3)-β-d-glucan (BG), and 16S rRNA gene sequencing to determine the circulating bacteriome.


Results
id="Par3">Of the 19 patients, 13 were in intensive care and 10 patients received mechanical ventilation.

 We found 8 patients with high EAA (≥ 0.6) and about half of the patients had high serum BG levels which tended to be higher in later in the illness.

 Although only 1 patient had a positive blood culture, 18 of 19 patients were positive for 16S rRNA gene amplification.

 Proteobacteria was the most abundant phylum.

 The diversity of bacterial genera was decreased overtime.


Conclusions
id="Par4">Bacterial DNA and toxins were discovered in virtually all severely ill COVID-19 pneumonia patients.

 This raises a previously unrecognized concern for significant contribution of bacterial products in the pathogenesis of this disease.



'@
$e12 = @'
[Phatadon%Sirivongrangson%NULL%1,  Win%Kulvichit%NULL%1,  Sunchai%Payungporn%NULL%1,  Trairak%Pisitkun%NULL%1,  Ariya%Chindamporn%NULL%1,  Sadudee%Peerapornratana%NULL%1,  Prapaporn%Pisitkun%NULL%1,  Suwalak%Chitcharoen%NULL%1,  Vorthon%Sawaswong%NULL%1,  Navaporn%Worasilchai%NULL%1,  Sarinya%Kampunya%NULL%1,  Opass%Putcharoen%NULL%1,  Thammasak%Thawitsri%NULL%1,  Nophol%Leelayuwatanakul%NULL%1,  Napplika%Kongpolprom%NULL%1,  Vorakamol%Phoophiboon%NULL%1,  Thitiwat%Sriprasart%NULL%1,  Rujipat%Samransamruajkit%NULL%1,  Somkanya%Tungsanga%NULL%1,  Kanitha%Tiankanon%NULL%1,  Nuttha%Lumlertgul%NULL%1,  Asada%Leelahavanichkul%NULL%1,  Tueboon%Sriphojanart%NULL%1,  Terapong%Tantawichien%NULL%1,  Usa%Thisyakorn%NULL%1,  Chintana%Chirathaworn%NULL%1,  Kearkiat%Praditpornsilpa%NULL%1,  Kriang%Tungsanga%NULL%1,  Somchai%Eiam-Ong%NULL%1,  Visith%Sitprija%NULL%1,  John A.%Kellum%NULL%1,  Nattachai%Srisawat%drnattachai@yahoo.com%1]
'@
$ws.Range("D12").Value = $d12
$ws.Range("E12").Value = $e12

$d20 = @'
Background
id="Par1">A COVID-19 outbreak started in Wuhan, China, last December and now has become a global pandemic.

 The clinical information in caring of critically ill patients with COVID-19 needs to be shared timely, especially under the situations that there is still a largely ongoing spread of COVID-19 in many countries.


Methods
id="Par2">A multicenter prospective observational study investigated all the COVID-19 patients received in 19 ICUs of 16 hospitals in Wuhan, China, over 24 h between 8 AM February 2h and 8 AM February 27, 2020. The demographic information, clinical characteristics, vital signs, complications, laboratory values, and clinical managements of the patients were studied.


Results
id="Par3">A total of 226 patients were included.

 Their median (interquartile range, IQR) age was 64 (57–70) years, and 139 (61.5%) patients were male.

 The duration from the date of ICU admission to the study date was 11 (5–17) days, and the duration from onset of symptoms to the study date was 31 (24–36) days.

 Among all the patients, 155 (68.6%) had at least one coexisting disease, and their sequential organ failure assessment score was 4 (2–8).

 Organ function damages were found in most of the patients: ARDS in 161 (71.2%) patients, septic shock in 34 (15.0%) patients, acute kidney injury occurred in 57 (25.2%) patients, cardiac injury in 61 (27.0%) patients, and lymphocytopenia in 160 (70.8%) patients.

 Of all the studied patients, 85 (37.6%) received invasive mechanical ventilation, including 14 (6.2%) treated with extracorporeal membrane oxygenation (ECMO) at the same time, 20 (8.8%) received noninvasive mechanical ventilation, and 24 (10.6%) received continuous renal replacement therapy.

 By April 9, 2020, 87 (38.5%) patients were deceased and 15 (6.7%) were still in the hospital.


Conclusions
id="Par4">Critically ill patients with COVID-19 are associated with a higher risk of severe complications and need to receive an intensive level of treatments.

 COVID-19 poses a great strain on critical care resources in hospitals.


Trial registration
id="Par5">Chinese Clinical Trial Registry, ChiCTR2000030164. Registered on February 24, 2020, http://www.

chictr.

org.

cn/edit.

aspx?pid=49983&amp;htm=4

'@
$e20 = @'
[Yuan%Yu%NULL%0,  Dan%Xu%NULL%0,  Shouzhi%Fu%NULL%0,  Jun%Zhang%NULL%0,  Xiaobo%Yang%NULL%0,  Liang%Xu%NULL%0,  Jiqian%Xu%NULL%0,  Yongran%Wu%NULL%0,  Chaolin%Huang%NULL%0,  Yaqi%Ouyang%NULL%0,  Luyu%Yang%NULL%0,  Minghao%Fang%NULL%0,  Hongwen%Xiao%NULL%0,  Jing%Ma%NULL%0,  Wei%Zhu%NULL%0,  Song%Hu%NULL%0,  Quan%Hu%NULL%0,  Daoyin%Ding%NULL%0,  Ming%Hu%NULL%0,  Guochao%Zhu%NULL%0,  Weijiang%Xu%NULL%0,  Jun%Guo%NULL%0,  Jinglong%Xu%NULL%0,  Haitao%Yuan%NULL%0,  Bin%Zhang%NULL%0,  Zhui%Yu%yuzhui@whu.edu.cn%0,  Dechang%Chen%icudechangchen@163.com%0,  Shiying%Yuan%yuan_shiying@163.com%0,  You%Shang%you_shanghust@163.com%0]
'@
$ws.Range("D20").Value = $d20
$ws.Range("E20").Value = $e20

$d25 = @'
Background
id="Par1">Coronavirus Disease 2019 (Covid-19) is expanding worldwide.

 The characteristics of this infection in patients varies from country to country.

 To move forward, clinical data on infected patients are needed.

 Here, we report a comparison between fatalities and recovery of patients with severe Covid-19, based on demographic and clinical characteristics.


Methods
id="Par2">Between 5 March and 12 May 2020 in Mashhad, Iran, 1278 of 4000 suspected Covid-19 patients were confirmed positive by real-time reverse-transcriptase–polymerase-chain-reaction assay of upper respiratory specimens.

 We compared the demographic, exposure history and clinical symptoms of 925 survivors and 353 fatal cases with confirmed disease.


Results
id="Par3">Mean (SD) age for all confirmed patients was 56.9 (18.7) years, 67.1 (15.9) years in fatal cases and 53.0 (18.3) years in survivors.

 Multivariate logistic regression analysis showed that the outcome of patients was associated with age (odds ratio = 1.049, P = 0.0001, 95% CI = 1.040–1.057).

 Despite a high burden of Covid-19 infections in the 30–39 and 40–49 year age groups, most of these (89.6 and 87.2%, respectively) recovered.

 The median (IQR) duration of hospitalization was 9.0 (6.0–14.0) days.

 The most prevalent co-morbidities were cardiovascular disorders (21%) and diabetes (16.3%).

 Dyspnoea (72.7%), cough (68.1%) and fever (63.8%) were the most frequent clinical symptoms.

 Healthcare workers, of whom two (3%) died, comprised 5.2% of infected cases.

 Combination antiviral and antibiotic therapy was used in 43.0% of cases.


Conclusions
id="Par4">The characteristics of severe Covid-19 varied substantially between fatal cases and survivors, with diabetes and cardiovascular disorders the most prevalent co-morbidities.

 In contrast to other studies, there were a higher number of fatalities in younger patients in our setting.


Supplementary Information
The online version contains supplementary material available at 10.1186/s12879-021-06363-6.

'@
$e25 = @'
[Ladan%Goshayeshi%NULL%1,  Mina%Akbari Rad%NULL%1,  Robert%Bergquist%NULL%1,  Abolghasem%Allahyari%NULL%1,  Kamila%Hashemzadeh%NULL%1,  Nasrin%Milani%NULL%1,  Mahdi%Gholian-Aval%NULL%1,  Fariba%Rezaeitalab%NULL%1,  Mohammad Jafar%Sadeghi Quchani%NULL%1,  Zahra%Nahbandani%NULL%1,  Mandana%Khodashahi%NULL%1,  Zahra%Javid%NULL%1,  Mahnaz%Mozdourian%NULL%1,  Mohammad Ali%Yaghoubi%NULL%1,  Zahra%Mozaheb%NULL%1,  Mohsen%Seddigh-Shamsi%NULL%1,  Mohammad%Moeini Nodeh%NULL%1,  Shima%Nabavi%NULL%1,  Hooman%Mosannen Mozaffari%NULL%1,  Mohammadreza%Farzanehfar%NULL%1,  Zahra%Lotfi%NULL%0,  Alireza%Shariati%NULL%1,  Shekoofe%Bonakdaran%NULL%1,  Zahra%Rezaieyazdi%NULL%1,  Zahra%Mirfeizi%NULL%1,  Maryam%Miri%NULL%1,  Reza%Bassiri%NULL%1,  Sajjad%Ataei Azimi%NULL%1,  Benyamin%Hoseini%Hoseinib1@nums.ac.ir%3,  Benyamin%Hoseini%Hoseinib1@nums.ac.ir%0,  Benyamin%Hoseini%Hoseinib1@nums.ac.ir%0]
'@
$ws.Range("D25").Value = $d25
$ws.Range("E25").Value = $e25

$d27 = @'
Background
id="Par1">In March 2020, the WHO declared the novel coronavirus (COVID-19) outbreak a global pandemic.

 Although the number of infected cases is increasing, information about its clinical characteristics in the Middle East, especially in Iran, a country which is considered to be one of the most important focal points of the disease in the world, is lacking.

 To date, there is no available literature on the clinical data on COVID-19 patients in Iran.


Methods
id="Par2">In this multicenter retrospective study, 113 hospitalized confirmed cases of COVID-19 admitted to university affiliated hospitals in Shiraz, Iran from February 20 to March 20 were entered in the study.


Results
id="Par3">The mean age was 53.75 years and 71 (62.8%) were males.

 The most common symptoms at onset were fatigue (75: 66.4%), cough (73: 64.6%), and fever (67: 59.3%).

 Laboratory data revealed significant correlation between lymphocyte count (P value = 0.003), partial thromboplastin time (P value = 0.000), international normalized ratio (P value = 0.000) with the severity of the disease.

 The most common abnormality in chest CT scans was ground-glass opacity (77: 93.9%), followed by consolidation (48: 58.5%).

 Our results revealed an overall 8% (9 out of 113 cases) mortality rate among patients, in which the majority was among patients admitted to the ICU (5: 55.6%).


Conclusion
id="Par4">Evaluating the clinical data of COVID-19 patients and finding the source of infection and studying the behavior of the disease is crucial for understanding the pandemic.



'@
$e27 = @'
[Reza%Shahriarirad%NULL%1,  Zohre%Khodamoradi%NULL%1,  Amirhossein%Erfani%NULL%1,  Hamidreza%Hosseinpour%NULL%1,  Keivan%Ranjbar%NULL%1,  Yasaman%Emami%NULL%1,  Alireza%Mirahmadizadeh%NULL%1,  Mehrzad%Lotfi%NULL%1,  Babak%Shirazi Yeganeh%NULL%1,  Abolfazl%Dorrani Nejad%NULL%1,  Abdolrasool%Hemmati%NULL%1,  Mostafa%Ebrahimi%NULL%1,  Mohsen%Moghadami%moghadami@sums.ac.ir%1]
'@
$ws.Range("D27").Value = $d27
$ws.Range("E27").Value = $e27

$d28 = @'
Background
id="Par1">Acute kidney injury (AKI) is a severe complication of coronavirus disease-2019 (COVID-19).

 This study aims to evaluate incidence, risk factors and case-fatality rate of AKI in patients with COVID-19.
Methods
id="Par2">We reviewed the health medical records of 307 consecutive patients with COVID-19 hospitalized at the University Hospital of Modena, Italy.


Results
id="Par3">AKI was diagnosed in 69 out of 307 (22.4%) COVID-19 patients.

 Stages 1, 2, or 3 AKI accounted for 57.9%, 24.6% and 17.3%, respectively.

 AKI patients had a mean age of 74.7 ± 9.9 years.

 These patients showed higher serum levels of the main markers of inflammation and higher rate of severe pneumonia than non-AKI patients.

 Kidney injury was associated with a higher rate of urinary abnormalities including proteinuria (0.44 ± 0.85 vs 0.18 ± 0.29 mg/mg; P =  &lt; 0.0001) and microscopic hematuria (P = 0.032) compared to non-AKI patients.

 Hemodialysis was performed in 7.2% of the subjects and 33.3% of the survivors did not recover kidney function after AKI.

 Risk factors for kidney injury were age, male sex, CKD and higher non-renal SOFA score.

 Patients with AKI had a mortality rate of 56.5%.

 Adjusted Cox regression analysis revealed that COVID-19-associated AKI was independently associated with in-hospital death (hazard ratio [HR] = 4.82; CI 95%, 1.36–17.08) compared to non-AKI patients.


Conclusion
id="Par4">AKI was a common and harmful consequence of COVID-19. It manifested with urinary abnormalities (proteinuria, microscopic hematuria) and conferred an increased risk for death.

 Given the well-known short-term sequelae of AKI, prevention of kidney injury is imperative in this vulnerable cohort of patients.


Supplementary Information
The online version contains supplementary material available at 10.1007/s10157-021-02092-x.



'@
$e28 = @'
[Gaetano%Alfano%gaetano.alfano@unimore.it%1,  Annachiara%Ferrari%NULL%2,  Annachiara%Ferrari%NULL%0,  Francesco%Fontana%NULL%1,  Giacomo%Mori%NULL%1,  Riccardo%Magistroni%NULL%1,  Marianna%Meschiari%NULL%1,  Erica%Franceschini%NULL%1,  Marianna%Menozzi%NULL%1,  Gianluca%Cuomo%NULL%1,  Gabriella%Orlando%NULL%1,  Antonella%Santoro%NULL%1,  Margherita%Digaetano%NULL%1,  Cinzia%Puzzolante%NULL%1,  Federica%Carli%NULL%1,  Andrea%Bedini%NULL%1,  Jovana%Milic%NULL%1,  Irene%Coloretti%NULL%1,  Paolo%Raggi%NULL%1,  Cristina%Mussini%NULL%1,  Massimo%Girardis%NULL%1,  Gianni%Cappelli%NULL%1,  Giovanni%Guaraldi%NULL%1,  NULL%NULL%NULL%0]
'@
$ws.Range("D28").Value = $d28
$ws.Range("E28").Value = $e28

$d30 = @'
id="Par1">Invasive mechanical ventilation (IMV) is the standard treatment in critically ill COVID-19 patients with acute severe respiratory distress syndrome (ARDS).

 When IMV setting is extremely aggressive, especially through the application of high positive-end-expiratory respiration (PEEP) values, lung damage can occur.

 Until today, in COVID-19 patients, two types of ARDS were identified (L- and H-type); for the L-type, a lower PEEP strategy was supposed to be preferred, but data are still missing.

 The aim of this study was to evaluate if a clinical management with lower PEEP values in critically ill L-type COVID-19 patients was safe and efficient in comparison to usual standard of care.

 A retrospective analysis was conducted on consecutive patients with COVID-19 ARDS admitted to the ICU and treated with IMV.

 Patients were treated with a lower PEEP strategy adapted to BMI: PEEP 10 cmH2O if BMI &lt; 30 kg m−2, PEEP 12 cmH2O if BMI 30–50 kg m−2, PEEP 15 cmH2O if BMI &gt; 50 kg m−2.
 Primary endpoint was the PaO2/FiO2 ratio evolution during the first 3 IMV days; secondary endpoints were to analyze ICU length of stay (LOS) and IMV length.

 From March 2 to January 15, 2021, 79 patients underwent IMV.

 Average applied PEEP was 11 ± 2.9 cmH2O for BMI &lt; 30 kg m−2 and 16 ± 3.18 cmH2O for BMI &gt; 30 kg m−2.
 During the first 24 h of IMV, patients’ PaO2/FiO2 ratio presented an improvement (p&lt;0.001; CI 99%) that continued daily up to 72 h (p&lt;0.001; CI 99%).

 Median ICU LOS was 15 days (10–28); median duration of IMV was 12 days (8–26).

 The ICU mortality rate was 31.6%.

 Lower PEEP strategy treatment in L-type COVID-19 ARDS resulted in a PaO2/FiO2 ratio persistent daily improvement during the first 72 h of IMV.

 A lower PEEP strategy could be beneficial in the first phase of ARDS in critically ill COVID-19 patients.


'@
$e30 = @'
[Samuele%Ceruti%samuele.ceruti@moncucco.ch%1,  Marco%Roncador%NULL%2,  Marco%Roncador%NULL%0,  Andrea%Saporito%NULL%1,  Maira%Biggiogero%NULL%1,  Andrea%Glotta%NULL%1,  Pier Andrea%Maida%NULL%1,  Patrizia%Urso%NULL%1,  Giovanni%Bona%NULL%1,  Christian%Garzoni%NULL%1,  Romano%Mauri%NULL%0,  Alain%Borgeat%NULL%1]
'@
$ws.Range("D30").Value = $d30
$ws.Range("E30").Value = $e30

$d33 = @'
id="Par1">Critically ill patients with COVID-19 pneumonia suffered both high thrombotic and bleeding risk.

 The effect of SARS-CoV-2 on coagulation and fibrinolysis is not well known.

 We conducted a retrospective study of critically ill patients admitted to an intensive care unit (ICU) a cause of severe COVID-19 pneumonia and we evaluated coagulation function using rotational thromboelastometry (ROTEM) on day of admission (T0) and 5 (T5) and 10 (T10) days after admission to ICU.

 Coagulation standard parameters were also evaluated.

 Forty patients were enrolled into the study.

 The ICU and the hospital mortality were 10% and 12.5%, respectively.

 On ICU admission, prothrombin time was slightly reduced and it increased significantly at T10 (T0 = 65.1 ± 9.8 vs T10 = 85.7 ± 1.5, p = 0.002), while activated partial thromboplastin time and fibrinogen values were higher at T0 than T10 (32.2 ± 2.9 vs 27.2 ± 2.1, p = 0.017 and 895.1 ± 110 vs 332.5 ± 50, p = 0.002, respectively); moreover, whole blood thromboelastometry profiles were consistent with hypercoagulability characterized by an acceleration of the propagation phase of blood clot formation [i.

e.

, CFT below the lower limit in INTEM 16/40 patients (40%) and EXTEM 20/40 patients (50%)] and significant higher clot strength [MCF above the upper limit in INTEM 20/40 patients (50%), in EXTEM 28/40 patients (70%) and in FIBTEM 29/40 patients (72.5%)]; however, this hypercoagulable state persists in the first five days, but it decreases ten day after, without returning to normal values.

 No sign of secondary hyperfibrinolysis or sepsis induced coagulopathy (SIC) were found during the study period.

 In six patients (15%) a deep vein thrombosis and in 2 patients (5%) a thromboembolic event, were found; 12 patients (30%) had a catheter-related thrombosis.

 ROTEM analysis confirms that patients with severe COVID-19 pneumonia had a hypercoagulation state that persisted over time.


'@
$e33 = @'
[Vittorio%Pavoni%NULL%1,  Lara%Gianesello%gianesello.lara@libero.it%1,  Maddalena%Pazzi%NULL%2,  Maddalena%Pazzi%NULL%0,  Caterina%Stera%NULL%1,  Tommaso%Meconi%NULL%1,  Francesca Covani%Frigieri%NULL%1]
'@
$ws.Range("D33").Value = $d33
$ws.Range("E33").Value = $e33

